$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (ECs, Pomc, Oprk1 -> ECs) ---
$ws.Cells.Item(2,4).Value  = "ECs"                     # D2
$ws.Cells.Item(2,5).Value  = 2                          # E2
$ws.Cells.Item(2,6).Value  = 0.6666666666666666         # F2
$ws.Cells.Item(2,7).Value  = 0.3971766666666667         # G2
$ws.Cells.Item(2,8).Value  = 1.19153                    # H2
$ws.Cells.Item(2,13).Value = 0.0005943333333333333      # M2
$ws.Cells.Item(2,14).Value = 0.001783                   # N2
$ws.Cells.Item(2,15).Value = 0.399059982094897          # O2
$ws.Cells.Item(2,16).Value = 0.3990599820948971         # P2
$ws.Cells.Item(2,17).Value = 0.0002360553322222222      # Q2
$ws.Cells.Item(2,18).Value = 0.00212449799              # R2
$ws.Cells.Item(2,19).Value = 0.399059982094897          # S2
$ws.Cells.Item(2,20).Value = 0.3990599820948971         # T2

# --- Add new row 3 (ECs, Pomc, Oprk1 -> FAPs) ---
$ws.Cells.Item(3,1).Value  = "ECs"                      # A3
$ws.Cells.Item(3,2).Value  = "Pomc"                     # B3
$ws.Cells.Item(3,3).Value  = "Oprk1"                    # C3
$ws.Cells.Item(3,4).Value  = "FAPs"                     # D3
$ws.Cells.Item(3,5).Value  = 2                          # E3
$ws.Cells.Item(3,6).Value  = 0.6666666666666666         # F3
$ws.Cells.Item(3,7).Value  = 0.3971766666666667         # G3
$ws.Cells.Item(3,8).Value  = 1.19153                    # H3
$ws.Cells.Item(3,9).Value  = 1                          # I3
$ws.Cells.Item(3,10).Value = 1                          # J3
$ws.Cells.Item(3,11).Value = 1                          # K3
$ws.Cells.Item(3,12).Value = 0.3333333333333333         # L3
$ws.Cells.Item(3,13).Value = 0.000895                   # M3
$ws.Cells.Item(3,14).Value = 0.002685                   # N3
$ws.Cells.Item(3,15).Value = 0.6009400179051029         # O3
$ws.Cells.Item(3,16).Value = 0.600940017905103          # P3
$ws.Cells.Item(3,17).Value = 0.0003554731166666666      # Q3
$ws.Cells.Item(3,18).Value = 0.00319925805              # R3
$ws.Cells.Item(3,19).Value = 0.6009400179051029         # S3
$ws.Cells.Item(3,20).Value = 0.600940017905103          # T3
